# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.7527432677738641;  E = 0.4942365360607697;  F = 0; G = 6.189590430959694 }
    3 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 0.1494219747398047;  E = 0.4942365360607697;  F = 0; G = 5.586269137925634 }
    4 = @{ B = 3.286832544864788;  C = 1.655778082260271;  D = 22.3905356188092;    E = 10.19245300693656;   F = 0; G = 37.52559925287081 }
    5 = @{ B = 0.01293466051926884; C = 0.002571899574220771; D = 3.537761648806719; E = 0.4942365360607697; F = 0; G = 4.047504744960978 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
